$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mahmoudreza Entezami is renaming their skill from "Writing Unit Testing"
# to "Automated Unit Testing" in column F (Skill), row 6.
$ws.Range("F6").Value = "Automated Unit Testing"

# Reflect the author's post-edit selection in the saved view state.
$ws.Range("F11").Select()
